$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the data set. It belongs right
# above the existing row 73 (chronologically it sits first in this run of
# entries), so insert a fresh row there — this pushes the former rows
# 73..110 down to 74..111, matching the rest of the sheet.
$ws.Rows("73:73").Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A73").Value2 = 10
$ws.Range("B73").Value2 = "Vega Modelo de Temuco"
$ws.Range("C73").Value2 = "La Araucanía"
$ws.Range("D73").Value2 = 44460
$ws.Range("E73").Value2 = 9
$ws.Range("F73").Value2 = 100114007
$ws.Range("G73").Value2 = "Jengibre"
$ws.Range("H73").Value2 = "Sin especificar"
$ws.Range("I73").Value2 = "Primera"
$ws.Range("J73").Value2 = 30
$ws.Range("K73").Value2 = 20000
$ws.Range("L73").Value2 = 20000
$ws.Range("M73").Value2 = 20000
$ws.Range("N73").Value2 = "`$/caja 13 kilos"
$ws.Range("O73").Value2 = "Perú"
$ws.Range("P73").Value2 = 1538
$ws.Range("Q73").Value2 = 13
$ws.Range("R73").Value2 = "Hortaliza"
